$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all data rows ("data"/"error" header + the 10 invalid-customer
# records) and any cell-level formatting (the bold header font + thin
# border + center/top alignment) that was applied to A1:B11.
# Clearing the whole used range collapses the sheet back down to the
# single default cell (A1) with no stored content, matching the
# dimension shrinking from A1:B11 to A1:A1 and sheetData becoming empty.
$ws.Cells.Clear()
